$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.673.31"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "2.616.94"
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'321.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "'109.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.538"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.09%  "
$ws.Range("D10").Value = "'39.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("D11").Value = "'19.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.41%  "
$ws.Range("D12").Value = "'0.0807"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "'7.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "3.037.54"
$ws.Range("E15").Value = "  +3.21%  "
$ws.Range("D16").Value = "2.595.67"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").Value = "'0.859"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "48.680.32"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").Value = "'2.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.04%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'12.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.28%  "
$ws.Range("D21").Value = "'6.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").Value = "'268.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.70%  "
$ws.Range("D24").Value = "'68.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.66%  "
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "'25.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'9.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("E30").Value = "  -4.28%  "
$ws.Range("D31").Value = "'34.65"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("D32").Value = "'49.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "'19.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("D36").Value = "'0.0791"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").Value = "'4.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.03%  "
$ws.Range("D38").Value = "'2.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").Value = "'3.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.29%  "
$ws.Range("D40").Value = "'125.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.92%  "
$ws.Range("D41").Value = "'22.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("E43").Value = "  -2.95%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").Value = "2.054.33"
$ws.Range("E45").Value = "  +2.15%  "
$ws.Range("D46").Value = "'3.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.45%  "
$ws.Range("D47").Value = "'2.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.23%  "
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").Value = "'8.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("D50").Value = "'58.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.31%  "
$ws.Range("D51").Value = "'5.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.11%  "
